$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.328.18'
$ws.Range('E2').Value = '  -0.61%  '
$ws.Range('D3').Value = '2.641.62'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '517.54'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.95'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.573'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '2.649.26'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').Value = '3.107.72'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').Value = '59.332.19'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.03'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('D18').Value = '2.636.08'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '348.78'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.31'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -2.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.21'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.51'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.165'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +2.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.996'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('D28').Value = '0.0₃0804'
$ws.Range('E28').Value = '  -2.55%  '
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.55'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.91'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.70'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.06'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.946'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -12.44%  '
$ws.Range('E37').Value = '  +0.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.863'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.65'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('E40').Value = '  +3.16%  '
$ws.Range('E41').Value = '  -1.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '278.19'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0990'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.62'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -1.52%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.600'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -3.25%  '
$ws.Range('D47').Value = '2.079.08'
$ws.Range('E47').Value = '  +4.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0529'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -3.21%  '
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('E51').Value = '  -0.58%  '
